$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-header the table: drop "ID EMPLEADO"/"FECHA HORA"/"DIRECCION"/
# "NOMBRE DISPOSITIVO"/"SERIAL DISPOSITIVO" columns and reorder/rename the
# remaining ones to match the payroll-calc layout. ---
$ws.Range("A1").Value = "HORARIO"
$ws.Range("B1").Value = "NOMBRE"
$ws.Range("C1").Value = "FECHA Y HORA"
$ws.Range("D1").Value = "FECHA"
$ws.Range("E1").Value = "HORA"
$ws.Range("F1").Value = "DIA"

# Remove the now-unused last two columns (G: NOMBRE DISPOSITIVO, H: SERIAL DISPOSITIVO)
$ws.Range("G1:H1").EntireColumn.Delete()

# Re-apply the AutoFilter so its range shrinks to the new A1:F1 extent.
# Toggling it off then on again forces Excel to recompute the filter range
# from the current region instead of keeping the stale A1:H1 reference.
$ws.Range("A1:F1").AutoFilter()
$ws.Range("A1:F1").AutoFilter()

# Sync the hidden _FilterDatabase defined name to the new range as well.
$names = $wb.Names
for ($i = 1; $i -le $names.Count; $i++) {
    $nm = $names.Item($i)
    if ($nm.Name -like "*_FilterDatabase*") {
        $nm.RefersTo = "=Asistencia!`$A`$1:`$F`$1"
    }
}

# Match the saved selection/active cell from the edit.
$ws.Range("F1").Select()
